# components.xlsx - "added nfc antenna components, added most of footprints"
#
# Adds a handful of new component rows (with price/qty/formula), fills in a
# couple of rows that previously only had a price/qty but no name/link, and
# fixes up one product link that used to point at a generic manufacturer
# page. Finishes by restoring the original selection/scroll position as
# closely as the object model allows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 9: FPC display connector ---------------------------------
$ws.Range("A9").Value = "FPC display conn"
$ws.Range("B9").Value = "https://www.mouser.ee/ProductDetail/GCT/FFC2B35-24-T?qs=Li%252BoUPsLEnvLbIbV0OhDVA%3D%3D"
$ws.Range("C9").Value = 0.47
$ws.Range("D9").Value = 1
$ws.Range("E9").Formula = "=C9*D9"

# --- Row 16: serial-to-parallel shift register (previously blank) -----
$ws.Range("A16").Value = "serial to parralel"
$ws.Range("B16").Value = "https://www.mouser.ee/ProductDetail/Nexperia/74HC595PW-Q100118?qs=1sbE9T7hb3aHrTORCcEuDg%3D%3D"
$ws.Range("C16").Value = 0.43
$ws.Range("D16").Value = 1

# --- Row 29: display schottky - replace generic manufacturer link -----
# with the specific product page (name/price/qty unchanged).
$ws.Range("B29").Value = "https://www.mouser.ee/ProductDetail/Toshiba/CUHS20S30H3F?qs=PqoDHHvF64%252BnIC9Qnnw9zg%3D%3D"

# --- Row 32: NFC antenna connector (JST XH 3p) -------------------------
$ws.Range("A32").Value = "JST XH 3p"
$ws.Range("B32").Value = "https://www.aliexpress.com/item/1005003559631954.html"

# --- Row 34: current-sense resistor footprint --------------------------
$ws.Range("A34").Value = "Resistor 5mO"
$ws.Range("B34").Value = "https://www.mouser.ee/ProductDetail/Vishay/WFCP06125L000FE66?qs=sGAEpiMZZMtlubZbdhIBIJBDgjsVQBFlUYEtiJor9t8%3D"
$ws.Range("C34").Value = 0.29
$ws.Range("D34").Value = 2

# --- Restore cursor / scroll position (best effort) --------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("C40").Select()
